{"js": "// ----------------------------------------------------------------------\n// Edit: Paytrack ERP integration doc \u2014 narrow scope down to \"Reembolso\"\n// (remove \"Adiantamentos em esp\u00e9cie\" scenario everywhere), rewrite the\n// main analysis paragraph, and drop the now-irrelevant field rows from\n// the mapping table.\n// ----------------------------------------------------------------------\n\n// 1) \"Processos desejados: Adiantamentos em esp\u00e9cie, Reembolso\"\n//    -> \"Processos desejados: Reembolso\"\nlet results = context.document.body.search(\"Adiantamentos em esp\u00e9cie, \", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) \"Informa\u00e7\u00f5es necess\u00e1rias pelo ERP: Valor total relat\u00f3rio, CPF, Empresa,\n//     Descri\u00e7\u00e3o relat\u00f3rio, Tipo de documento, Valor do rateio, Centro de\n//     custo, Ordem interna, Conta cont\u00e1bil\"\n//    -> \"Informa\u00e7\u00f5es necess\u00e1rias pelo ERP: Valor total relat\u00f3rio, CPF\"\nresults = context.document.body.search(\n  \", Empresa, Descri\u00e7\u00e3o relat\u00f3rio, Tipo de documento, Valor do rateio, Centro de custo, Ordem interna, Conta cont\u00e1bil\",\n  { matchCase: true }\n);\nresults.load(\"text\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Rewrite the big analysis paragraph (the one that still starts with\n//    \"Para realizar a integra\u00e7\u00e3o do ERP SAP ECC/4HANA\").\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Para realizar a integra\u00e7\u00e3o do ERP SAP ECC/4HANA\";\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(marker) === 0) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!targetParagraph) {\n  throw new Error(\"Could not locate the analysis paragraph to rewrite.\");\n}\n\nconst NEW_ANALYSIS_TEXT = \"Para realizar a integra\u00e7\u00e3o do ERP SAP ECC/4HANA com o seu SaaS Paytrack para o processo de reembolso, \u00e9 importante seguir as diretrizes fornecidas e criar uma an\u00e1lise funcional detalhada para cada cen\u00e1rio desejado. Vou te orientar sobre como estruturar essa an\u00e1lise funcional:\\u000b\\u000b**1. Identifica\u00e7\u00e3o do Cen\u00e1rio:** \\u000b   - Cen\u00e1rio: Reembolso\\u000b\\u000b**2. Informa\u00e7\u00f5es Necess\u00e1rias pelo ERP (SAP ECC/4HANA):**\\u000b   - Valor total do relat\u00f3rio\\u000b   - CPF\\u000b\\u000b**3. Mapeamento de Campos:**\\u000b\\u000b| Campo Paytrack | Campo SAP ECC/4HANA |\\u000b|-----------------|----------------------|\\u000b| Valor Total     | Betrg                |\\u000b| CPF             | PersNumber           |\\u000b\\u000b**4. JSON de Exemplo Formatado:**\\u000b\\u000b```json\\u000b{\\u000b  \\\"Reembolso\\\": {\\u000b    \\\"Valor Total\\\": \\\"1000.00\\\",\\u000b    \\\"CPF\\\": \\\"123.456.789-00\\\"\\u000b  }\\u000b}\\u000b```\\u000b\\u000b**5. Observa\u00e7\u00f5es Importantes:**\\u000b- Utiliza\u00e7\u00e3o de comunica\u00e7\u00e3o s\u00edncrona com os Webservices do cliente.\\u000b- A Paytrack ser\u00e1 ativa nas integra\u00e7\u00f5es, aguardando o cliente disponibilizar um Webservice para consumo.\\u000b- Separar a an\u00e1lise funcional por cen\u00e1rio selecionado: adiantamento, presta\u00e7\u00e3o de contas, etc.\\u000b\\u000bCom essa estrutura, voc\u00ea ter\u00e1 uma an\u00e1lise funcional clara e organizada para guiar a integra\u00e7\u00e3o do seu SaaS Paytrack com o ERP SAP ECC/4HANA no processo de reembolso. Certifique-se de documentar cada passo e manter uma boa comunica\u00e7\u00e3o com o cliente para garantir o sucesso da integra\u00e7\u00e3o.\";\n\nconst wholeRange = targetParagraph.getRange(\"Whole\");\nwholeRange.insertText(NEW_ANALYSIS_TEXT, Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) Drop the mapping-table rows that are no longer relevant: keep the\n//    header row, \"Valor total relat\u00f3rio\" and \"CPF\"; remove \"Empresa\",\n//    \"Descri\u00e7\u00e3o relat\u00f3rio\", \"Tipo de documento\", \"Valor do rateio\",\n//    \"Centro de custo\", \"Ordem interna\", \"Conta cont\u00e1bil\".\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  table.rows.load(\"items\");\n  await context.sync();\n\n  const rows = table.rows.items;\n  for (let i = 0; i < rows.length; i++) {\n    rows[i].cells.load(\"items\");\n  }\n  await context.sync();\n\n  for (let i = 0; i < rows.length; i++) {\n    const cells = rows[i].cells.items;\n    for (let j = 0; j < cells.length; j++) {\n      cells[j].body.load(\"text\");\n    }\n  }\n  await context.sync();\n\n  const keep = new Set([\"Campo\", \"Valor total relat\u00f3rio\", \"CPF\"]);\n  const indicesToDelete = [];\n  for (let i = 0; i < rows.length; i++) {\n    const firstCellText = rows[i].cells.items[0].body.text.trim();\n    if (!keep.has(firstCellText)) {\n      indicesToDelete.push(i);\n    }\n  }\n  // Delete from the highest index down, re-resolving the row object from\n  // the table each time \u2014 deleting a row shifts every later index, and a\n  // previously-loaded row proxy can become stale once its neighbors move.\n  for (let k = indicesToDelete.length - 1; k >= 0; k--) {\n    table.rows.items[indicesToDelete[k]].delete();\n    await context.sync();\n  }\n}\n", "ps1": "# ----------------------------------------------------------------------\n# Edit: Paytrack ERP integration doc - narrow scope down to \"Reembolso\"\n# (remove \"Adiantamentos em esp\u00e9cie\" scenario everywhere), rewrite the\n# main analysis paragraph, and drop the now-irrelevant field rows from\n# the mapping table.\n# ----------------------------------------------------------------------\n\n$d = $word.ActiveDocument\n\n# 1) \"Processos desejados: Adiantamentos em esp\u00e9cie, Reembolso\"\n#    -> \"Processos desejados: Reembolso\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Adiantamentos em esp\u00e9cie, \"\n$find.Replacement.Text = \"\"\n$find.Execute(\"Adiantamentos em esp\u00e9cie, \", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\n# 2) \"Informa\u00e7\u00f5es necess\u00e1rias pelo ERP: Valor total relat\u00f3rio, CPF, Empresa,\n#     Descri\u00e7\u00e3o relat\u00f3rio, Tipo de documento, Valor do rateio, Centro de\n#     custo, Ordem interna, Conta cont\u00e1bil\"\n#    -> \"Informa\u00e7\u00f5es necess\u00e1rias pelo ERP: Valor total relat\u00f3rio, CPF\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$removed = \", Empresa, Descri\u00e7\u00e3o relat\u00f3rio, Tipo de documento, Valor do rateio, Centro de custo, Ordem interna, Conta cont\u00e1bil\"\n$find2.Execute($removed, $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\n# 3) Rewrite the big analysis paragraph (the one that still starts with\n#    \"Para realizar a integra\u00e7\u00e3o do ERP SAP ECC/4HANA\").\n$marker = \"Para realizar a integra\u00e7\u00e3o do ERP SAP ECC/4HANA\"\n$targetParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.Length -ge $marker.Length -and $t.Substring(0, $marker.Length) -eq $marker) {\n        $targetParagraph = $p\n        break\n    }\n}\nif ($null -eq $targetParagraph) {\n    throw \"Could not locate the analysis paragraph to rewrite.\"\n}\n\n$NEW_ANALYSIS_TEXT = 'Para realizar a integra\u00e7\u00e3o do ERP SAP ECC/4HANA com o seu SaaS Paytrack para o processo de reembolso, \u00e9 importante seguir as diretrizes fornecidas e criar uma an\u00e1lise funcional detalhada para cada cen\u00e1rio desejado. Vou te orientar sobre como estruturar essa an\u00e1lise funcional:' + [char]11 + '' + [char]11 + '**1. Identifica\u00e7\u00e3o do Cen\u00e1rio:** ' + [char]11 + '   - Cen\u00e1rio: Reembolso' + [char]11 + '' + [char]11 + '**2. Informa\u00e7\u00f5es Necess\u00e1rias pelo ERP (SAP ECC/4HANA):**' + [char]11 + '   - Valor total do relat\u00f3rio' + [char]11 + '   - CPF' + [char]11 + '' + [char]11 + '**3. Mapeamento de Campos:**' + [char]11 + '' + [char]11 + '| Campo Paytrack | Campo SAP ECC/4HANA |' + [char]11 + '|-----------------|----------------------|' + [char]11 + '| Valor Total     | Betrg                |' + [char]11 + '| CPF             | PersNumber           |' + [char]11 + '' + [char]11 + '**4. JSON de Exemplo Formatado:**' + [char]11 + '' + [char]11 + '```json' + [char]11 + '{' + [char]11 + '  \"Reembolso\": {' + [char]11 + '    \"Valor Total\": \"1000.00\",' + [char]11 + '    \"CPF\": \"123.456.789-00\"' + [char]11 + '  }' + [char]11 + '}' + [char]11 + '```' + [char]11 + '' + [char]11 + '**5. Observa\u00e7\u00f5es Importantes:**' + [char]11 + '- Utiliza\u00e7\u00e3o de comunica\u00e7\u00e3o s\u00edncrona com os Webservices do cliente.' + [char]11 + '- A Paytrack ser\u00e1 ativa nas integra\u00e7\u00f5es, aguardando o cliente disponibilizar um Webservice para consumo.' + [char]11 + '- Separar a an\u00e1lise funcional por cen\u00e1rio selecionado: adiantamento, presta\u00e7\u00e3o de contas, etc.' + [char]11 + '' + [char]11 + 'Com essa estrutura, voc\u00ea ter\u00e1 uma an\u00e1lise funcional clara e organizada para guiar a integra\u00e7\u00e3o do seu SaaS Paytrack com o ERP SAP ECC/4HANA no processo de reembolso. Certifique-se de documentar cada passo e manter uma boa comunica\u00e7\u00e3o com o cliente para garantir o sucesso da integra\u00e7\u00e3o.'\n\n$rng = $targetParagraph.Range\n# Exclude the trailing paragraph mark so we only replace the paragraph's\n# visible content (keeps the paragraph itself, just swaps its runs/text).\n$rng.End = $rng.End - 1\n$rng.Text = $NEW_ANALYSIS_TEXT\n\n# 4) Drop the mapping-table rows that are no longer relevant: keep the\n#    header row, \"Valor total relat\u00f3rio\" and \"CPF\"; remove \"Empresa\",\n#    \"Descri\u00e7\u00e3o relat\u00f3rio\", \"Tipo de documento\", \"Valor do rateio\",\n#    \"Centro de custo\", \"Ordem interna\", \"Conta cont\u00e1bil\".\n$table = $d.Tables.Item(1)\n$keep = @(\"Campo\", \"Valor total relat\u00f3rio\", \"CPF\")\nfor ($i = $table.Rows.Count; $i -ge 1; $i--) {\n    $cellText = $table.Cell($i, 1).Range.Text.TrimEnd([char]7, [char]13)\n    if ($keep -notcontains $cellText) {\n        $table.Rows.Item($i).Delete()\n    }\n}\n"}
